$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price record for "Achicoria" was reported for the
# week of 2023-04-21 (serial 45037). Insert it as the new row 18 (the
# sheet is ordered with the most recent entries first among this block),
# pushing the existing rows 18-84 down to 19-85 and carrying the former
# last row (old row 84) out to the new row 85.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 45037
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100112010
$ws.Range("G18").Value = "Achicoria"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 10000
$ws.Range("N18").Value = "`$/caja 18 unidades"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 556
$ws.Range("Q18").Value = 18
$ws.Range("R18").Value = "Hortaliza"
